$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 2
$ws.Range("A2").Value = 2187390
$ws.Range("B2").Value = 440000
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 5

# Update the active selection on the sheet view
$ws.Range("G9").Select()
